$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 108.30769
$ws.Range("I11").Value = 108.30769
$ws.Range("K11").Value = 108.30769
$ws.Range("M11").Value = 31.69231000000001

$ws.Range("H17").Value = 8335500
$ws.Range("J17").Value = 8335500
$ws.Range("L17").Value = 25006500
$ws.Range("N17").Value = -25006836

$ws.Range("H33").Value = 19272.941
$ws.Range("I33").Value = 25078.309
$ws.Range("K33").Value = 25078.309
$ws.Range("M33").Value = -24849.309

$ws.Range("H92").Value = 381.4
$ws.Range("I92").Value = 436.69232
$ws.Range("K92").Value = 436.69232
$ws.Range("M92").Value = 811.30768

$ws.Range("H98").Value = 4705.622
$ws.Range("I98").Value = 4839.775
$ws.Range("J98").Value = 3632.4
$ws.Range("K98").Value = 4839.775
$ws.Range("L98").Value = 3632.4
$ws.Range("M98").Value = -3341.775
$ws.Range("N98").Value = -6628.4

$ws.Range("H100").Value = 1398.6
$ws.Range("I100").Value = 1420.6666
$ws.Range("K100").Value = 1420.6666
$ws.Range("M100").Value = -879.6666

$ws.Range("H122").Value = 4705.622
$ws.Range("I122").Value = 4839.775
$ws.Range("J122").Value = 3632.4
$ws.Range("K122").Value = 14519.325
$ws.Range("L122").Value = 10897.2
$ws.Range("M122").Value = -12069.325
$ws.Range("N122").Value = -15797.2

$ws.Range("H125").Value = 14026.917
$ws.Range("J125").Value = 9298.833000000001
$ws.Range("L125").Value = 83689.497
$ws.Range("N125").Value = -88609.497

$ws.Range("H132").Value = 3716.361
$ws.Range("I132").Value = 3953.8484
$ws.Range("K132").Value = 11861.5452
$ws.Range("M132").Value = -9331.5452

$ws.Range("H137").Value = 1572.4445
$ws.Range("I137").Value = 1425.25
$ws.Range("J137").Value = 2750
$ws.Range("K137").Value = 4275.75
$ws.Range("L137").Value = 8250
$ws.Range("M137").Value = -1725.75
$ws.Range("N137").Value = -13350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3816.4167
$ws.Range("I45").Value = 2690.0312
$ws.Range("J45").Value = 6069.1875
$ws.Range("K45").Value = 2690.0312
$ws.Range("L45").Value = 6069.1875
$ws.Range("M45").Value = -2313.0312
$ws.Range("N45").Value = -6823.1875

$ws.Range("H74").Value = 8326.77
$ws.Range("I74").Value = 6749.8
$ws.Range("J74").Value = 9312.375
$ws.Range("K74").Value = 6749.8
$ws.Range("L74").Value = 9312.375
$ws.Range("M74").Value = -5875.8
$ws.Range("N74").Value = -11060.375

$ws.Range("H77").Value = 8326.77
$ws.Range("I77").Value = 6749.8
$ws.Range("J77").Value = 9312.375
$ws.Range("K77").Value = 33749
$ws.Range("L77").Value = 46561.875
$ws.Range("M77").Value = -29381
$ws.Range("N77").Value = -55297.875

$ws.Range("H97").Value = 400.20834
$ws.Range("I97").Value = 285.9524
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 285.9524
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = 210.0476
$ws.Range("N97").Value = -2192

$ws.Range("H110").Value = 9287.5
$ws.Range("I110").Value = 7060
$ws.Range("J110").Value = 13000
$ws.Range("K110").Value = 7060
$ws.Range("L110").Value = 13000
$ws.Range("M110").Value = -5015
$ws.Range("N110").Value = -17090

$ws.Range("H132").Value = 4209.66
$ws.Range("I132").Value = 3400.6924
$ws.Range("K132").Value = 10202.0772
$ws.Range("M132").Value = -7672.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 252.66667
$ws.Range("J64").Value = 286.16666
$ws.Range("L64").Value = 286.16666
$ws.Range("N64").Value = -736.16666

$ws.Range("H67").Value = 252.66667
$ws.Range("J67").Value = 286.16666
$ws.Range("L67").Value = 286.16666
$ws.Range("N67").Value = -1846.16666

$ws.Range("H80").Value = 589.6667
$ws.Range("J80").Value = 738.625
$ws.Range("L80").Value = 738.625
$ws.Range("N80").Value = -2734.625

$ws.Range("H82").Value = 24284.143
$ws.Range("J82").Value = 39997.25
$ws.Range("L82").Value = 39997.25
$ws.Range("N82").Value = -40763.25

$ws.Range("H83").Value = 589.6667
$ws.Range("J83").Value = 738.625
$ws.Range("L83").Value = 3693.125
$ws.Range("N83").Value = -13677.125

$ws.Range("H85").Value = 24284.143
$ws.Range("J85").Value = 39997.25
$ws.Range("L85").Value = 39997.25
$ws.Range("N85").Value = -42649.25

$ws.Range("H86").Value = 100001176
$ws.Range("I86").Value = 2350
$ws.Range("K86").Value = 2350
$ws.Range("M86").Value = -1227

$ws.Range("H89").Value = 100001176
$ws.Range("I89").Value = 2350
$ws.Range("K89").Value = 11750
$ws.Range("M89").Value = -6134

$ws.Range("H94").Value = 2911.0625
$ws.Range("I94").Value = 2282.4167
$ws.Range("J94").Value = 4797
$ws.Range("K94").Value = 2282.4167
$ws.Range("L94").Value = 4797
$ws.Range("M94").Value = -1831.4167
$ws.Range("N94").Value = -5699

$ws.Range("H99").Value = 4803.4644
$ws.Range("I99").Value = 3779.85
$ws.Range("K99").Value = 3779.85
$ws.Range("M99").Value = -2281.85

$ws.Range("H107").Value = 2424.524
$ws.Range("I107").Value = 2003.7931
$ws.Range("J107").Value = 3363.077
$ws.Range("K107").Value = 2003.7931
$ws.Range("L107").Value = 3363.077
$ws.Range("M107").Value = -83.79310000000009
$ws.Range("N107").Value = -7203.077

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5560.696
$ws.Range("I31").Value = 4369.7856
$ws.Range("J31").Value = 7413.222
$ws.Range("K31").Value = 4369.7856
$ws.Range("L31").Value = 7413.222
$ws.Range("M31").Value = -4074.7856
$ws.Range("N31").Value = -8003.222

$ws.Range("H34").Value = 5560.696
$ws.Range("I34").Value = 4369.7856
$ws.Range("J34").Value = 7413.222
$ws.Range("K34").Value = 4369.7856
$ws.Range("L34").Value = 7413.222
$ws.Range("M34").Value = -4167.7856
$ws.Range("N34").Value = -7817.222

$ws.Range("H38").Value = 4268.5
$ws.Range("J38").Value = 4268.5
$ws.Range("L38").Value = 4268.5
$ws.Range("N38").Value = -5022.5

$ws.Range("H41").Value = 11900
$ws.Range("I41").Value = 7120
$ws.Range("J41").Value = 23850
$ws.Range("K41").Value = 7120
$ws.Range("L41").Value = 23850
$ws.Range("M41").Value = -6692
$ws.Range("N41").Value = -24706

$ws.Range("H46").Value = 4268.5
$ws.Range("J46").Value = 4268.5
$ws.Range("L46").Value = 4268.5
$ws.Range("N46").Value = -4690.5

$ws.Range("H59").Value = 40750
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H62").Value = 7502
$ws.Range("I62").Value = 6005
$ws.Range("K62").Value = 6005
$ws.Range("M62").Value = -5381

$ws.Range("H65").Value = 7502
$ws.Range("I65").Value = 6005
$ws.Range("K65").Value = 30025
$ws.Range("M65").Value = -26905

$ws.Range("H122").Value = 3994.85
$ws.Range("I122").Value = 3999.8
$ws.Range("K122").Value = 11999.4
$ws.Range("M122").Value = -9549.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 246.46153
$ws.Range("J12").Value = 439.42856
$ws.Range("L12").Value = 1318.28568
$ws.Range("N12").Value = -1664.28568

$ws.Range("H132").Value = 100000980
$ws.Range("J132").Value = 1333
$ws.Range("L132").Value = 11997
$ws.Range("N132").Value = -17057

$ws.Range("H137").Value = 33591.25
$ws.Range("I137").Value = 1352.5
$ws.Range("J137").Value = 65830
$ws.Range("K137").Value = 4057.5
$ws.Range("L137").Value = 197490
$ws.Range("M137").Value = 1042.5
$ws.Range("N137").Value = -207690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2084.889
$ws.Range("I97").Value = 623.8889
$ws.Range("K97").Value = 623.8889
$ws.Range("M97").Value = -127.8889

$ws.Range("H126").Value = 2419.1333
$ws.Range("I126").Value = 1403.8182
$ws.Range("K126").Value = 4211.4546
$ws.Range("M126").Value = -1741.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 1000000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H55").Value = 818.1875
$ws.Range("I55").Value = 831.1818
$ws.Range("K55").Value = 831.1818
$ws.Range("M55").Value = -658.1818

$ws.Range("H93").Value = 12031.381
$ws.Range("I93").Value = 1288.8462
$ws.Range("K93").Value = 1288.8462
$ws.Range("M93").Value = -40.84619999999995

$ws.Range("H136").Value = 8249
$ws.Range("I136").Value = 8249
$ws.Range("K136").Value = 24747
$ws.Range("M136").Value = -22197

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2941
$ws.Range("I81").Value = 1550
$ws.Range("K81").Value = 3100
$ws.Range("M81").Value = -2039

$ws.Range("H84").Value = 2941
$ws.Range("I84").Value = 1550
$ws.Range("K84").Value = 15500
$ws.Range("M84").Value = -10196

$ws.Range("H113").Value = 305.73914
$ws.Range("I113").Value = 326.21054
$ws.Range("J113").Value = 208.5
$ws.Range("K113").Value = 978.6316199999999
$ws.Range("L113").Value = 625.5
$ws.Range("M113").Value = 1191.36838
$ws.Range("N113").Value = -4965.5

$ws.Range("H136").Value = 8835.214
$ws.Range("I136").Value = 7569.3
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 22707.9
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -20157.9
$ws.Range("N136").Value = -41100
